# "UI cambiada y arreglados fallitos"
#  - Fix a content typo: the verb "LLORAR" in C2 should read "IR".
#  - Move the saved selection/active cell to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix cell C2: "LLORAR" -> "IR"
$ws.Range("C2").Value = "IR"

# Leave the selection on C6, matching the refreshed sheet view
$ws.Range("C6").Select()
